$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" (sheet1): insert a new row 3 for the handed-back
#     file b3fc4203-...md, pushing the old row 3 (f346e95e-...md) to row 4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Rows.Item(3).Insert()

$wsOverview.Range("A3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.md"
$wsOverview.Range("B3").Value = "e2e\b3fc4203-d0de-47ab-a55b-c5af32d89420.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-26 08:48:53"

# Resize the Overview table (table3) to cover the new row
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# Rebuild the hyperlinks for this sheet in final order
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bb62e1dfadc09451ba7548c290f07354364f4fe/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md", "", "", "e2e\97cdfdf1-cd49-46ea-b470-466da18b27ac.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4c76a1cddabfbfa9a4dc22976081ea6f0e9d591/e2e/b3fc4203-d0de-47ab-a55b-c5af32d89420.md", "", "", "e2e\b3fc4203-d0de-47ab-a55b-c5af32d89420.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de9911703c32f44f81919edf8ab951293bab3631/e2e/f346e95e-57e6-4aca-bea4-ce0ef35182e2.md", "", "", "e2e\f346e95e-57e6-4aca-bea4-ce0ef35182e2.md")

# --- Sheet "zh-cn" (sheet2): same new-row insert, with the zh-cn xliff columns ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Rows.Item(3).Insert()

$wsZh.Range("A3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "'Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("G3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.3cc3982bc12e3e9c7c44147e9f7cb39121e8bd7c.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-26 08:48:48"
$wsZh.Range("I3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.md"
$wsZh.Range("J3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.3cc3982bc12e3e9c7c44147e9f7cb39121e8bd7c.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-26 08:49:19"
$wsZh.Range("L3").Value = "'"
$wsZh.Range("M3").Value = "'True"
$wsZh.Range("N3").Value = "'"
$wsZh.Range("O3").Value = "'False"
$wsZh.Range("P3").Value = "'"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bb62e1dfadc09451ba7548c290f07354364f4fe/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md", "", "", "97cdfdf1-cd49-46ea-b470-466da18b27ac.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/c63c03a4989bad93d06f83ea4111ec6004325b93/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md", "", "", "97cdfdf1-cd49-46ea-b470-466da18b27ac.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4c76a1cddabfbfa9a4dc22976081ea6f0e9d591/e2e/b3fc4203-d0de-47ab-a55b-c5af32d89420.md", "", "", "b3fc4203-d0de-47ab-a55b-c5af32d89420.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b9a9ea58b19d4763bd25a06d7a3ffee1dd18c0c3/e2e/b3fc4203-d0de-47ab-a55b-c5af32d89420.md", "", "", "b3fc4203-d0de-47ab-a55b-c5af32d89420.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de9911703c32f44f81919edf8ab951293bab3631/e2e/f346e95e-57e6-4aca-bea4-ce0ef35182e2.md", "", "", "f346e95e-57e6-4aca-bea4-ce0ef35182e2.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/310846e8b930d4889d0a991f2b93fd826b1fa856/e2e/f346e95e-57e6-4aca-bea4-ce0ef35182e2.md", "", "", "f346e95e-57e6-4aca-bea4-ce0ef35182e2.md")

# --- Sheet "de-de" (sheet3): same new-row insert, with the de-de xliff columns ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Rows.Item(3).Insert()

$wsDe.Range("A3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "'Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("G3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.3cc3982bc12e3e9c7c44147e9f7cb39121e8bd7c.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-26 08:43:57"
$wsDe.Range("I3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.md"
$wsDe.Range("J3").Value = "b3fc4203-d0de-47ab-a55b-c5af32d89420.3cc3982bc12e3e9c7c44147e9f7cb39121e8bd7c.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-26 08:49:26"
$wsDe.Range("L3").Value = "'"
$wsDe.Range("M3").Value = "'True"
$wsDe.Range("N3").Value = "'"
$wsDe.Range("O3").Value = "'False"
$wsDe.Range("P3").Value = "'"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6bb62e1dfadc09451ba7548c290f07354364f4fe/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md", "", "", "97cdfdf1-cd49-46ea-b470-466da18b27ac.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/93a3178e20b656cf39f83c4bf8d95d21485e4708/e2e/97cdfdf1-cd49-46ea-b470-466da18b27ac.md", "", "", "97cdfdf1-cd49-46ea-b470-466da18b27ac.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f4c76a1cddabfbfa9a4dc22976081ea6f0e9d591/e2e/b3fc4203-d0de-47ab-a55b-c5af32d89420.md", "", "", "b3fc4203-d0de-47ab-a55b-c5af32d89420.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/7c0df6d2bf61f7e8fe7c52be9d2da7c5e42a5f73/e2e/b3fc4203-d0de-47ab-a55b-c5af32d89420.md", "", "", "b3fc4203-d0de-47ab-a55b-c5af32d89420.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de9911703c32f44f81919edf8ab951293bab3631/e2e/f346e95e-57e6-4aca-bea4-ce0ef35182e2.md", "", "", "f346e95e-57e6-4aca-bea4-ce0ef35182e2.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/de0dd2bfec2af57729c67fcda2536fea7db20b70/e2e/f346e95e-57e6-4aca-bea4-ce0ef35182e2.md", "", "", "f346e95e-57e6-4aca-bea4-ce0ef35182e2.md")
